$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Proceso Expulsado"/"Proceso Expulsor" columns of the task-expulsion
# table (J2:L6) now also show the priority of the process next to its name.
$ws.Range("J3").Value = "Proceso Expulsado (prioridad)"
$ws.Range("K3").Value = "Proceso Expulsor (prioridad)"

$ws.Range("K4").Value = "Riesgos (40)"
$ws.Range("K5").Value = "Riesgos (40)"
$ws.Range("K6").Value = "Volante (20)"

# Select the widened columns (J:L) and make them wider so the longer
# "(prioridad)" labels fit comfortably.
$ws.Columns("J:L").Select() | Out-Null
$ws.Columns("J:L").ColumnWidth = 28.7

"Task expulsion table updated"
